# 3_1_Journal-Andrei-PiresDonose.xlsx
# Add two new journal entries (documentation protocols finished + testing
# with Luuk) to the weekly journal sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Block "Documentation sur les protocols de test." (week of 46037 / 2026-01-15) ---
# The block header (row 47) already exists; fill in the first detail line (row 48).
$ws.Range("A48").Value = 46037
$ws.Range("B48").Value = "Documentation sur les protocols de test."
$ws.Range("D48").Value = 3.5

# --- New block reflection text (row 54), describing the day's work ---
$ws.Range("B54").Value = "Aujourd'hui, j'ai finalisé les protocoles de test. J'ai pris le temps d'optimiser la structure pour éliminer les redondances et simplifier la lecture. Le résultat est désormais plus condensé et efficace. Globalement, la journée a été productive."
# Row grows taller to fit the wrapped paragraph, matching Excel's auto-fit result.
$ws.Rows(54).RowHeight = 47.25

# --- Block "Tester l'application avec Luuk" (week of 46038 / 2026-01-16) ---
$ws.Range("A55").Value = 46038
$ws.Range("B55").Value = "Tester l'application avec Luuk"
$ws.Range("D55").Value = 0.5

# Additional time logged under the same task (no description on this line).
$ws.Range("D56").Value = 3

# Leave the cursor where the author ended up editing.
$ws.Activate()
$ws.Range("D57").Select()
